$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aperturas_Siniestros")

# --- Header cells (G1:H1) -------------------------------------------------
# Copy header style from A1 (bold/fill header style) onto G1:H1 first, then
# set the text so the style carries over without clobbering the new values.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$ws.Range("G1").Value() = "tipo_indexacion_severidad"
$ws.Range("H1").Value() = "medida_indexacion_severidad"

# --- Column G (tipo_indexacion_severidad) ----------------------------------
$ws.Range("G2").Value() = "Ninguna"
$ws.Range("G3").Value() = "Ninguna"
$ws.Range("G4").Value() = "Ninguna"
$ws.Range("G5").Value() = "Ninguna"
$ws.Range("G6").Value() = "Por fecha de ocurrencia"
$ws.Range("G7").Value() = "Por fecha de movimiento"
$ws.Range("G8").Value() = "Ninguna"
$ws.Range("G9").Value() = "Ninguna"

# --- Column H (medida_indexacion_severidad) ---------------------------------
$ws.Range("H2").Value() = "Ninguna"
$ws.Range("H3").Value() = "Ninguna"
$ws.Range("H4").Value() = "Ninguna"
$ws.Range("H5").Value() = "Ninguna"
$ws.Range("H7").Value() = "IPC"
$ws.Range("H6").Value() = "SMMLV"
$ws.Range("H8").Value() = "Ninguna"
$ws.Range("H9").Value() = "Ninguna"

# --- Column widths for the new columns -------------------------------------
$ws.Columns.Item(7).ColumnWidth() = 22.57
$ws.Columns.Item(8).ColumnWidth() = 25.7

# --- Data validation list on the new columns --------------------------------
# Target sqref is "G2:G9 H8:H9" (full column G, plus the bottom two rows of H,
# since H2:H7 are driven by the indexation measure rather than a free list).
$full = $ws.Range("G2:H9")
$full.Validation.Add(3, 1, 1, """Ninguna,Por fecha de ocurrencia,Por fecha de movimiento""")
$exclude = $ws.Range("H2:H7")
$exclude.Validation.Delete()

# --- Selection tweak (matches the saved cursor position in the workbook) ---
$ws.Activate()
$ws.Range("F4").Select()
